$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to be stored as TEXT (not
# auto-converted to a number/date by Excel), then drop the now-unneeded
# "@" number-format override so the cell's style stays untouched (matches
# the original file, where these data cells carry no explicit style).
function Set-TextCell([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# --- Row 2 (Bitcoin) ---
Set-TextCell "D2" "65.545.16"
Set-TextCell "E2" "  +1.10%  "

# --- Row 3 (Ethereum) ---
Set-TextCell "D3" "3.394.02"
Set-TextCell "E3" "  +0.13%  "

# --- Row 4 (TetherUSD) ---
Set-TextCell "E4" "  -0.10%  "

# --- Row 5 (BNB) ---
Set-TextCell "D5" "560.92"
Set-TextCell "E5" "  +0.13%  "

# --- Row 6 (Solana) ---
Set-TextCell "D6" "176.02"
Set-TextCell "E6" "  +0.43%  "

# --- Row 7 (XRP) ---
Set-TextCell "D7" "0.632"
Set-TextCell "E7" "  +0.74%  "

# --- Row 8 (LidoStakedEther) ---
Set-TextCell "D8" "3.382.21"
Set-TextCell "E8" "  +0.11%  "

# --- Row 9 (USDC) ---
Set-TextCell "E9" "  -0.12%  "

# --- Row 10 (Dogecoin) ---
Set-TextCell "D10" "0.173"
Set-TextCell "E10" "  +4.72%  "

# --- Row 11 (Cardano) ---
Set-TextCell "E11" "  +0.69%  "

# --- Row 12 (Avalanche) ---
Set-TextCell "D12" "53.43"
Set-TextCell "E12" "  -2.04%  "

# --- Row 13 (ShibaInu) ---
Set-TextCell "E13" "  +0.22%  "

# --- Row 14 (Polkadot) ---
Set-TextCell "D14" "9.22"
Set-TextCell "E14" "  +0.80%  "

# --- Row 15 (WrappedliquidstakedEther2.0) ---
Set-TextCell "D15" "3.935.38"
Set-TextCell "E15" "  -0.49%  "

# --- Row 16 (Chainlink) ---
Set-TextCell "D16" "18.30"
Set-TextCell "E16" "  +0.07%  "

# --- Row 17 (WrappedEther) ---
Set-TextCell "D17" "3.400.27"
Set-TextCell "E17" "  -0.05%  "

# --- Row 18 (TRON) ---
Set-TextCell "E18" "  +1.03%  "

# --- Row 19 (WrappedBTC) ---
Set-TextCell "D19" "65.413.67"
Set-TextCell "E19" "  +0.86%  "

# --- Row 20 (Uniswap) ---
Set-TextCell "D20" "11.85"
Set-TextCell "E20" "  -0.54%  "

# --- Row 21 (Polygon) ---
Set-TextCell "E21" "  +0.48%  "

# --- Row 22 (BitcoinCash) ---
Set-TextCell "D22" "481.67"
Set-TextCell "E22" "  +2.50%  "

# --- Row 23 (Toncoin) ---
Set-TextCell "E23" "  -0.94%  "

# --- Row 24 / 25 swap: Litecoin <-> InternetComputer(DFINITY) ---
Set-TextCell "B24" "InternetComputer(DFINITY)"
Set-TextCell "C24" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D24" "14.35"
Set-TextCell "E24" "  +4.39%  "

Set-TextCell "B25" "Litecoin"
Set-TextCell "C25" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D25" "89.55"
Set-TextCell "E25" "  +3.52%  "

# --- Row 26 (PancakeSwap) ---
Set-TextCell "D26" "4.11"
Set-TextCell "E26" "  -0.49%  "

# --- Row 27 (ImmutableX) ---
Set-TextCell "E27" "  +1.72%  "

# --- Row 28 (RenderToken) ---
Set-TextCell "D28" "10.64"
Set-TextCell "E28" "  -2.05%  "

# --- Row 29 (Filecoin) ---
Set-TextCell "D29" "8.74"
Set-TextCell "E29" "  -1.17%  "

# --- Row 30 (EthereumClassic) ---
Set-TextCell "D30" "31.29"
Set-TextCell "E30" "  +2.06%  "

# --- Row 31 (NEARProtocol) ---
Set-TextCell "D31" "6.56"
Set-TextCell "E31" "  -2.78%  "

# --- Row 32 (Cosmos) ---
Set-TextCell "D32" "11.51"
Set-TextCell "E32" "  -0.27%  "

# --- Row 33 (OKB) ---
Set-TextCell "D33" "62.88"
Set-TextCell "E33" "  +4.85%  "

# --- Row 34 (Bittensor) ---
Set-TextCell "D34" "575.34"
Set-TextCell "E34" "  -0.71%  "

# --- Row 35 (Hedera) ---
Set-TextCell "E35" "  -0.88%  "

# --- Row 36 (Dai) ---
Set-TextCell "E36" "  +0.10%  "

# --- Row 37 (Stacks) ---
Set-TextCell "D37" "3.64"
Set-TextCell "E37" "  +5.09%  "

# --- Row 38 (Kaspa) ---
Set-TextCell "E38" "  +0.49%  "

# --- Row 39 (InjectiveProtocol) ---
Set-TextCell "D39" "35.85"
Set-TextCell "E39" "  -0.26%  "

# --- Row 40 (TheGraph) ---
Set-TextCell "E40" "  +0.37%  "

# --- Row 41 (PEPE) ---
Set-TextCell "D41" "0.0₃0740"
Set-TextCell "E41" "  -1.96%  "

# --- Row 42 (Maker) ---
Set-TextCell "D42" "3.098.40"
Set-TextCell "E42" "  -0.35%  "

# --- Row 43 (ThetaToken) ---
Set-TextCell "E43" "  -2.43%  "

# --- Row 44 (VeChain) ---
Set-TextCell "E44" "  +0.93%  "

# --- Row 45 (Stellar) ---
Set-TextCell "D45" "0.135"
Set-TextCell "E45" "  +0.28%  "

# --- Row 46 (ApeXProtocol) ---
Set-TextCell "D46" "3.16"
Set-TextCell "E46" "  -1.42%  "

# --- Row 47 (Fetch.AI) ---
Set-TextCell "E47" "  -3.67%  "

# --- Row 48 (FirstDigitalUSD) ---
Set-TextCell "E48" "  +0.00%  "

# --- Row 49 (Monero) ---
Set-TextCell "D49" "140.26"
Set-TextCell "E49" "  +2.71%  "

# --- Row 50 (WEMIXToken) ---
Set-TextCell "E50" "  -0.04%  "

# --- Row 51 (THORChain) ---
Set-TextCell "D51" "8.43"
Set-TextCell "E51" "  +0.56%  "
